$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.893987000000001
$ws.Range("H2").Value = 11.681961
$ws.Range("I2").Value = 0.143037189732266
$ws.Range("J2").Value = 0.143037189732266
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.319497666666667
$ws.Range("N2").Value = 27.958493
$ws.Range("O2").Value = 0.4629603708903243
$ws.Range("P2").Value = 0.4629603708903243
$ws.Range("Q2").Value = 36.29000276053034
$ws.Range("R2").Value = 326.610024844773
$ws.Range("S2").Value = 0.06622055040955956
$ws.Range("T2").Value = 0.06622055040955956

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.893987000000001
$ws.Range("H3").Value = 11.681961
$ws.Range("I3").Value = 0.143037189732266
$ws.Range("J3").Value = 0.143037189732266
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.09547066666667
$ws.Range("N3").Value = 30.286412
$ws.Range("O3").Value = 0.5015080223550378
$ws.Range("P3").Value = 0.5015080223550378
$ws.Range("Q3").Value = 39.31163153488134
$ws.Range("R3").Value = 353.8046838139321
$ws.Range("S3").Value = 0.07173429814585104
$ws.Range("T3").Value = 0.07173429814585104

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.893987000000001
$ws.Range("H4").Value = 11.681961
$ws.Range("I4").Value = 0.143037189732266
$ws.Range("J4").Value = 0.143037189732266
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.7152593333333334
$ws.Range("N4").Value = 2.145778
$ws.Range("O4").Value = 0.03553160675463796
$ws.Range("P4").Value = 0.03553160675463796
$ws.Range("Q4").Value = 2.785210545628667
$ws.Range("R4").Value = 25.066894910658
$ws.Range("S4").Value = 0.005082341176855415
$ws.Range("T4").Value = 0.005082341176855415

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.9333873333333332
$ws.Range("H5").Value = 2.800162
$ws.Range("I5").Value = 0.03428596476867894
$ws.Range("J5").Value = 0.03428596476867894
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.319497666666667
$ws.Range("N5").Value = 27.958493
$ws.Range("O5").Value = 0.4629603708903243
$ws.Range("P5").Value = 0.4629603708903243
$ws.Range("Q5").Value = 8.698701075096222
$ws.Range("R5").Value = 78.288309675866
$ws.Range("S5").Value = 0.01587304296564019
$ws.Range("T5").Value = 0.01587304296564019

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.9333873333333332
$ws.Range("H6").Value = 2.800162
$ws.Range("I6").Value = 0.03428596476867894
$ws.Range("J6").Value = 0.03428596476867894
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.09547066666667
$ws.Range("N6").Value = 30.286412
$ws.Range("O6").Value = 0.5015080223550378
$ws.Range("P6").Value = 0.5015080223550378
$ws.Range("Q6").Value = 9.422984444304889
$ws.Range("R6").Value = 84.806859998744
$ws.Range("S6").Value = 0.01719468638567467
$ws.Range("T6").Value = 0.01719468638567467

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.9333873333333332
$ws.Range("H7").Value = 2.800162
$ws.Range("I7").Value = 0.03428596476867894
$ws.Range("J7").Value = 0.03428596476867894
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.7152593333333334
$ws.Range("N7").Value = 2.145778
$ws.Range("O7").Value = 0.03553160675463796
$ws.Range("P7").Value = 0.03553160675463796
$ws.Range("Q7").Value = 0.6676140017817778
$ws.Range("R7").Value = 6.008526016036
$ws.Range("S7").Value = 0.001218235417364072
$ws.Range("T7").Value = 0.001218235417364072

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.286864
$ws.Range("H8").Value = 18.860592
$ws.Range("I8").Value = 0.2309343505227297
$ws.Range("J8").Value = 0.2309343505227297
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 9.319497666666667
$ws.Range("N8").Value = 27.958493
$ws.Range("O8").Value = 0.4629603708903243
$ws.Range("P8").Value = 0.4629603708903243
$ws.Range("Q8").Value = 58.59041437865067
$ws.Range("R8").Value = 527.313729407856
$ws.Range("S8").Value = 0.1069134525693191
$ws.Range("T8").Value = 0.1069134525693191

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.286864
$ws.Range("H9").Value = 18.860592
$ws.Range("I9").Value = 0.2309343505227297
$ws.Range("J9").Value = 0.2309343505227297
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 10.09547066666667
$ws.Range("N9").Value = 30.286412
$ws.Range("O9").Value = 0.5015080223550378
$ws.Range("P9").Value = 0.5015080223550378
$ws.Range("Q9").Value = 63.46885109732268
$ws.Range("R9").Value = 571.2196598759041
$ws.Range("S9").Value = 0.1158154294244993
$ws.Range("T9").Value = 0.1158154294244993

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.286864
$ws.Range("H10").Value = 18.860592
$ws.Range("I10").Value = 0.2309343505227297
$ws.Range("J10").Value = 0.2309343505227297
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.7152593333333334
$ws.Range("N10").Value = 2.145778
$ws.Range("O10").Value = 0.03553160675463796
$ws.Range("P10").Value = 0.03553160675463796
$ws.Range("Q10").Value = 4.496738153397334
$ws.Range("R10").Value = 40.470643380576
$ws.Range("S10").Value = 0.008205468528911354
$ws.Range("T10").Value = 0.008205468528911354

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 8.572307333333333
$ws.Range("H11").Value = 25.716922
$ws.Range("I11").Value = 0.3148851679477345
$ws.Range("J11").Value = 0.3148851679477346
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.319497666666667
$ws.Range("N11").Value = 27.958493
$ws.Range("O11").Value = 0.4629603708903243
$ws.Range("P11").Value = 0.4629603708903243
$ws.Range("Q11").Value = 79.88959819094956
$ws.Range("R11").Value = 719.006383718546
$ws.Range("S11").Value = 0.1457793541409452
$ws.Range("T11").Value = 0.1457793541409453

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 8.572307333333333
$ws.Range("H12").Value = 25.716922
$ws.Range("I12").Value = 0.3148851679477345
$ws.Range("J12").Value = 0.3148851679477346
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.09547066666667
$ws.Range("N12").Value = 30.286412
$ws.Range("O12").Value = 0.5015080223550378
$ws.Range("P12").Value = 0.5015080223550378
$ws.Range("Q12").Value = 86.54147722931822
$ws.Range("R12").Value = 778.873295063864
$ws.Range("S12").Value = 0.1579174378464023
$ws.Range("T12").Value = 0.1579174378464023

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 8.572307333333333
$ws.Range("H13").Value = 25.716922
$ws.Range("I13").Value = 0.3148851679477345
$ws.Range("J13").Value = 0.3148851679477346
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.7152593333333334
$ws.Range("N13").Value = 2.145778
$ws.Range("O13").Value = 0.03553160675463796
$ws.Range("P13").Value = 0.03553160675463796
$ws.Range("Q13").Value = 6.131422828368445
$ws.Range("R13").Value = 55.182805455316
$ws.Range("S13").Value = 0.01118837596038703
$ws.Range("T13").Value = 0.01118837596038704

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.116573333333334
$ws.Range("H14").Value = 12.34972
$ws.Range("I14").Value = 0.1512134172319493
$ws.Range("J14").Value = 0.1512134172319493
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 9.319497666666667
$ws.Range("N14").Value = 27.958493
$ws.Range("O14").Value = 0.4629603708903243
$ws.Range("P14").Value = 0.4629603708903243
$ws.Range("Q14").Value = 38.36439557466223
$ws.Range("R14").Value = 345.27956017196
$ws.Range("S14").Value = 0.07000581972529661
$ws.Range("T14").Value = 0.07000581972529661

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.116573333333334
$ws.Range("H15").Value = 12.34972
$ws.Range("I15").Value = 0.1512134172319493
$ws.Range("J15").Value = 0.1512134172319493
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 10.09547066666667
$ws.Range("N15").Value = 30.286412
$ws.Range("O15").Value = 0.5015080223550378
$ws.Range("P15").Value = 0.5015080223550378
$ws.Range("Q15").Value = 41.5587453338489
$ws.Range("R15").Value = 374.0287080046401
$ws.Range("S15").Value = 0.0758347418295421
$ws.Range("T15").Value = 0.0758347418295421

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.116573333333334
$ws.Range("H16").Value = 12.34972
$ws.Range("I16").Value = 0.1512134172319493
$ws.Range("J16").Value = 0.1512134172319493
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.7152593333333334
$ws.Range("N16").Value = 2.145778
$ws.Range("O16").Value = 0.03553160675463796
$ws.Range("P16").Value = 0.03553160675463796
$ws.Range("Q16").Value = 2.944417498017779
$ws.Range("R16").Value = 26.49975748216
$ws.Range("S16").Value = 0.005372855677110619
$ws.Range("T16").Value = 0.005372855677110619

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.420479333333333
$ws.Range("H17").Value = 10.261438
$ws.Range("I17").Value = 0.1256439097966415
$ws.Range("J17").Value = 0.1256439097966415
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 9.319497666666667
$ws.Range("N17").Value = 27.958493
$ws.Range("O17").Value = 0.4629603708903243
$ws.Range("P17").Value = 0.4629603708903243
$ws.Range("Q17").Value = 31.87714916588155
$ws.Range("R17").Value = 286.894342492934
$ws.Range("S17").Value = 0.05816815107956359
$ws.Range("T17").Value = 0.05816815107956359

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 3.420479333333333
$ws.Range("H18").Value = 10.261438
$ws.Range("I18").Value = 0.1256439097966415
$ws.Range("J18").Value = 0.1256439097966415
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 10.09547066666667
$ws.Range("N18").Value = 30.286412
$ws.Range("O18").Value = 0.5015080223550378
$ws.Range("P18").Value = 0.5015080223550378
$ws.Range("Q18").Value = 34.53134877560623
$ws.Range("R18").Value = 310.782138980456
$ws.Range("S18").Value = 0.06301142872306843
$ws.Range("T18").Value = 0.06301142872306843

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 3.420479333333333
$ws.Range("H19").Value = 10.261438
$ws.Range("I19").Value = 0.1256439097966415
$ws.Range("J19").Value = 0.1256439097966415
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 0.7152593333333334
$ws.Range("N19").Value = 2.145778
$ws.Range("O19").Value = 0.03553160675463796
$ws.Range("P19").Value = 0.03553160675463796
$ws.Range("Q19").Value = 2.446529767640444
$ws.Range("R19").Value = 22.018767908764
$ws.Range("S19").Value = 0.00446432999400947
$ws.Range("T19").Value = 0.00446432999400947
